$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows of incense log data to append below the existing table (rows 24-35)
# Date values are Excel serial date numbers (matching the existing column A data)
$data = @(
    @(45463, "Night",     "Goregasm",      1),
    @(45469, "Night",     "Goregasm",      1),
    @(45474, "Night",     "Nag Champa",    1),
    @(45477, "Night",     "Coffin Candy",  1),
    @(45481, "Morning",   "Dragons Blood", 1),
    @(45485, "Afternoon", "Goregasm",      1),
    @(45488, "Night",     "Goregasm",      1),
    @(45491, "Afternoon", "Dragons Blood", 1),
    @(45494, "Night",     "Nightshade",    1),
    @(45501, "Night",     "Nightshade",    1),
    @(45505, "Morning",   "Dragons Blood", 1),
    @(45507, "Night",     "Goregasm",      1)
)

$startRow = 24

# Copy the date formatting (style) from the last existing date cell (A23) so the
# new date cells reuse the same style instead of creating a new one.
$ws.Cells.Item($startRow - 1, 1).Copy()

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $row[0]
    $cellA.PasteSpecial(-4122) # xlPasteFormats

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$ws.Range("A36").Select()
